$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, pushing existing rows 5-12 down to 6-13
$ws.Rows.Item(5).Insert()

# Fill in the new row 5 with the latest weekly price entry
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C5").Value = "Los Lagos"
$ws.Range("D5").Value = 44481
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 100112013
$ws.Range("G5").Value = "Alcachofa"
$ws.Range("H5").Value = "Madrigal"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("N5").Value = "$/caja 50 unidades"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 220
$ws.Range("Q5").Value = 50
$ws.Range("R5").Value = "Hortaliza"
